$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info (card holder name / card number) ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking card number that must stay TEXT (as in
# the source file) rather than being auto-converted to a number. Format
# the cell as Text first, write the digits, then re-apply the original
# cell format (copied from a sibling cell) so the visual style/formatting
# index is unchanged.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 29.07.2025"

# --- Transaction row 6 ---
$ws.Range("B6").Value = "01.08."
$ws.Range("C6").Value = "02.08."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 88290831"
$ws.Range("E6").Value = "84,71-"

# --- Transaction row 7 ---
$ws.Range("B7").Value = "04.08."
$ws.Range("C7").Value = "05.08."
$ws.Range("D7").Value = "PAYPAL FJFKBU"
$ws.Range("E7").Value = "75,27-"

# --- Transaction row 8 ---
$ws.Range("B8").Value = "07.08."
$ws.Range("C8").Value = "08.08."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-72141532"
$ws.Range("E8").Value = "54,03-"

# --- Rows 9-11 no longer have transactions this cycle: clear them out ---
$ws.Range("B9:D9").ClearContents()
$ws.Range("E9:F9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

$ws.Range("B10:D10").ClearContents()
$ws.Range("E10:F10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

$ws.Range("B11:D11").ClearContents()
$ws.Range("E11:F11").Value = ""
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 12.08.2025"
$ws.Range("E12").Value = "214,01-"

# --- Next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.08.2025"
